$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transaksi")

$ws.Range("A3").Value = "Header Booking Order"
$ws.Range("B3").Value = "TRX010"
$ws.Range("A4").Value = "Detail Booking Order"
$ws.Range("B4").Value = "TRX011"
$ws.Range("A5").Value = "Detail Pembayaran -  Booking Order"
$ws.Range("B5").Value = "TRX012"

$ws.Columns.Item(1).ColumnWidth = 32.33

$ws.Range("B10").Select() | Out-Null
